$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 21 de Marzo de 2020 a las 15:46"

# Rows 29-32 resort: Tenerife's case count grew (192 -> 219), pushing it above
# Murcia/Cantabria/Leon in the "Casos totales" descending sort.
$ws.Range("A29").Value = "Tenerife"
$ws.Range("B29").Value = 219
$ws.Range("C29").Value = 4
$ws.Range("D29").Value = 211
$ws.Range("E29").Value = 4

$ws.Range("A30").Value = "Murcia"
$ws.Range("B30").Value = 215
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 213
$ws.Range("E30").Value = 1

$ws.Range("A31").Value = "Cantabria"
$ws.Range("B31").Value = 215
$ws.Range("C31").Value = 11
$ws.Range("D31").Value = 200
$ws.Range("E31").Value = 4

$ws.Range("A32").Value = "Leon"
$ws.Range("B32").Value = 201
$ws.Range("C32").Value = 3
$ws.Range("D32").Value = 156
$ws.Range("E32").Value = 7

# Rows 43-45 resort: Gran Canaria's case count grew (70 -> 97), pushing it
# above Cuenca/Ourense in the same sort.
$ws.Range("A43").Value = "Gran Canaria"
$ws.Range("B43").Value = 97
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 96
$ws.Range("E43").Value = 1

$ws.Range("A44").Value = "Cuenca"
$ws.Range("B44").Value = 94
$ws.Range("C44").Value = 5
$ws.Range("D44").Value = 84
$ws.Range("E44").Value = 5

$ws.Range("A45").Value = "Ourense"
$ws.Range("B45").Value = 74
$ws.Range("C45").Value = 5
$ws.Range("D45").Value = 74
$ws.Range("E45").Value = 0

# Simple numeric updates, no reordering.
$ws.Range("B55").Value = 13
$ws.Range("D55").Value = 13

$ws.Range("B56").Value = 11
$ws.Range("D56").Value = 11

$ws.Range("B59").Value = 4
$ws.Range("D59").Value = 4
